# MAJ automatique BRVM via GitHub Actions
# Refresh the daily "recommandations" sheet: updated day-count / variation
# figures for each title, plus the knock-on re-sort of several tied rows
# (same values shuffle to a new row position as rankings move).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Sheet1')

$ws.Range("C2").Value = 128
$ws.Range("D2").Value = 56241.74
$ws.Range("E2").Value = 120.2
$ws.Range("C3").Value = 51
$ws.Range("D3").Value = 46285
$ws.Range("E3").Value = 950
$ws.Range("C4").Value = 64
$ws.Range("D4").Value = 43571.11
$ws.Range("E4").Value = 667.6
$ws.Range("C5").Value = 52
$ws.Range("D5").Value = 43175
$ws.Range("A6").Value = "NEI-CEDA CI"
$ws.Range("C6").Value = 59
$ws.Range("D6").Value = 42620
$ws.Range("E6").Value = 730
$ws.Range("A7").Value = "SAFCA CI"
$ws.Range("C7").Value = 53
$ws.Range("D7").Value = 42085
$ws.Range("E7").Value = 800
$ws.Range("C8").Value = 64
$ws.Range("D8").Value = 37885
$ws.Range("E8").Value = 605
$ws.Range("C9").Value = 64
$ws.Range("D9").Value = 36930
$ws.Range("E9").Value = 550
$ws.Range("C10").Value = 64
$ws.Range("D10").Value = 31750
$ws.Range("E10").Value = 505
$ws.Range("C11").Value = 64
$ws.Range("D11").Value = 25655
$ws.Range("E11").Value = 395
$ws.Range("C12").Value = 64
$ws.Range("D12").Value = 23542.77
$ws.Range("E12").Value = 380.11
$ws.Range("C13").Value = 64
$ws.Range("D13").Value = 21136.33
$ws.Range("E13").Value = 338.12
$ws.Range("C14").Value = 64
$ws.Range("D14").Value = 13737.75
$ws.Range("E14").Value = 233.8
$ws.Range("C15").Value = 64
$ws.Range("D15").Value = 9400.35
$ws.Range("E15").Value = 156.46
$ws.Range("C16").Value = 64
$ws.Range("D16").Value = 8416.41
$ws.Range("E16").Value = 139.28
$ws.Range("C17").Value = 64
$ws.Range("D17").Value = 7746.69
$ws.Range("E17").Value = 121.93
$ws.Range("C18").Value = 64
$ws.Range("D18").Value = 7206.16
$ws.Range("E18").Value = 119.66
$ws.Range("C19").Value = 64
$ws.Range("D19").Value = 7135.18
$ws.Range("E19").Value = 119.07
$ws.Range("C20").Value = 64
$ws.Range("D20").Value = 6980.1
$ws.Range("E20").Value = 116.54
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 6881.77
$ws.Range("E21").Value = 113.83
$ws.Range("C22").Value = 64
$ws.Range("D22").Value = 6763.2
$ws.Range("E22").Value = 111.87
$ws.Range("C23").Value = 64
$ws.Range("D23").Value = 6419.97
$ws.Range("E23").Value = 95.11
$ws.Range("C24").Value = 64
$ws.Range("D24").Value = 6358.52
$ws.Range("E24").Value = 97.06
$ws.Range("A37").Value = "SAFCA CI (SAFC)"
$ws.Range("B37").Value = 8
$ws.Range("C37").Value = 4
$ws.Range("D37").Value = 22.3
$ws.Range("E37").Value = 1.27
$ws.Range("A38").Value = "PALM CI (PALC)"
$ws.Range("C38").Value = 5
$ws.Range("D38").Value = 21.55
$ws.Range("E38").Value = -4.17
$ws.Range("A39").Value = "BERNABE CI (BNBC)"
$ws.Range("B39").Value = 13
$ws.Range("C39").Value = 9
$ws.Range("D39").Value = 20.16
$ws.Range("E39").Value = 6
$ws.Range("A40").Value = "CFAO MOTORS CI (CFAC)"
$ws.Range("B40").Value = 7
$ws.Range("C40").Value = 8
$ws.Range("D40").Value = 16.48
$ws.Range("E40").Value = -6.56
$ws.Range("A43").Value = "BANK OF AFRICA NG (BOAN)"
$ws.Range("B43").Value = 13
$ws.Range("C43").Value = 12
$ws.Range("D43").Value = 13.55
$ws.Range("E43").Value = 4.49
$ws.Range("A44").Value = "SUCRIVOIRE (SCRC)"
$ws.Range("B44").Value = 8
$ws.Range("C44").Value = 8
$ws.Range("D44").Value = 12.76
$ws.Range("E44").Value = -1.01
$ws.Range("C45").Value = 6
$ws.Range("D45").Value = 11.06
$ws.Range("E45").Value = -1.28
$ws.Range("C48").Value = 6
$ws.Range("D48").Value = 7.68
$ws.Range("E48").Value = -1.83
$ws.Range("A51").Value = "VIVO ENERGY CI (SHEC)"
$ws.Range("B51").Value = 4
$ws.Range("C51").Value = 4
$ws.Range("D51").Value = 2.92
$ws.Range("E51").Value = 2.63
$ws.Range("A52").Value = "UNILEVER CI (UNLC)"
$ws.Range("B52").Value = 7
$ws.Range("C52").Value = 7
$ws.Range("D52").Value = 2.53
$ws.Range("E52").Value = -7.46
$ws.Range("A53").Value = "SETAO CI (STAC)"
$ws.Range("B53").Value = 13
$ws.Range("C53").Value = 11
$ws.Range("D53").Value = 2.27
$ws.Range("E53").Value = 1.83
$ws.Range("C54").Value = 63
$ws.Range("A57").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws.Range("B57").Value = 9
$ws.Range("C57").Value = 14
$ws.Range("D57").Value = -6.43
$ws.Range("E57").Value = 2.34
$ws.Range("A58").Value = "SOGB CI (SOGC)"
$ws.Range("B58").Value = 5
$ws.Range("C58").Value = 5
$ws.Range("D58").Value = -6.54
$ws.Range("E58").Value = 2.78
$ws.Range("A59").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws.Range("B59").Value = 6
$ws.Range("C59").Value = 9
$ws.Range("D59").Value = -6.92
$ws.Range("E59").Value = -1.35
$ws.Range("A60").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws.Range("B60").Value = 0
$ws.Range("D60").Value = -9.130000000000001
$ws.Range("E60").Value = -2.63
$ws.Range("A62").Value = "SOLIBRA CI (SLBC)"
$ws.Range("B62").Value = 10
$ws.Range("C62").Value = 13
$ws.Range("D62").Value = -12.23
$ws.Range("E62").Value = 3.91
$ws.Range("A63").Value = "ORAGROUP TOGO (ORGT)"
$ws.Range("B63").Value = 6
$ws.Range("C63").Value = 8
$ws.Range("D63").Value = -12.37
$ws.Range("E63").Value = 0.31
$ws.Range("A64").Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws.Range("B64").Value = 3
$ws.Range("C64").Value = 8
$ws.Range("D64").Value = -17.86
$ws.Range("E64").Value = -2.78
